# LoginData.xlsx update: add a second "Sheet2" worksheet with an extra
# set of login credentials (username/pass headers + two credential rows),
# the first credential row mirroring Sheet1's hyperlinked e-mail/password,
# and make Sheet2 the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- add the new worksheet right after Sheet1 -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# --- header row -------------------------------------------------------------
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "pass"

# --- credential rows ---------------------------------------------------------
$ws2.Range("A2").Value = "arun.joseph@learnship.com"
$ws2.Range("B2").Value = "Airtel@123"
$ws2.Range("A3").Value = "Insightadmin"
$ws2.Range("B3").Value = "Insight_0217"

# --- hyperlinks on the first credential row (e-mail + password got auto-linked) ---
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:arun.joseph@learnship.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:Airtel@123") | Out-Null

# A3 carries the same "Hyperlink" cell style (leftover formatting) but no live link
$ws2.Range("A3").Style = "Hyperlink"

# --- size the columns to fit their contents ----------------------------------
$ws2.Columns("A:B").AutoFit() | Out-Null

# --- make Sheet2 the active sheet/selection -----------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("D4").Select() | Out-Null

Write-Output "done"
